$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force text format to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.519.68"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.618.54"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.34"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.524"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.83"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("D9").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.847.77"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.616.28"
$ws.Range("D13").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.548"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.01"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.494.06"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.30"
$ws.Range("D18").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("D20").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.99"
$ws.Range("D25").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("D29").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.445.32"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.06"
$ws.Range("D34").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.938"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.562"
$ws.Range("D38").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.80"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.992"
$ws.Range("D43").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.758.45"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.70"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.41"
$ws.Range("D49").ClearFormats()

# Column E (Volume/1h) updates - plain text, safe to assign directly
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  +6.11%  "
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  +5.00%  "
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("E45").Value = "  -4.60%  "
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  +13.14%  "
$ws.Range("E51").Value = "  +1.64%  "
